$wb = $excel.ActiveWorkbook

# 1. Asthma Status: new row 8, A8 = 3 + 1
$ws1 = $wb.Worksheets.Item("Asthma Status")
$ws1.Range("A8").Formula = "=3  + 1"

# 2. Demographics: new row 13, A13 = 8 + 'Asthma Status'!A8
$ws2 = $wb.Worksheets.Item("Demographics")
$ws2.Range("A13").Formula = "= 8 +'Asthma Status'!A8"

# 3. Exercise: new row 7, A7 = 2 + Demographics!A13
$ws3 = $wb.Worksheets.Item("Exercise")
$ws3.Range("A7").Formula = "= 2 +Demographics!A13"

# 4. Health Status: new row 11, A11 = 6 + Exercise!A7
$ws4 = $wb.Worksheets.Item("Health Status")
$ws4.Range("A11").Formula = "= 6 +Exercise!A7"

# 5. Tobbaco USe: new row 17, A17 = 12 + 'Health Status'!A11
$ws5 = $wb.Worksheets.Item("Tobbaco USe")
$ws5.Range("A17").Formula = "= 12 +'Health Status'!A11"

# 6. Alcohol Consumption: new row 12, A12 = 7 + 'Tobbaco USe'!A17
$ws6 = $wb.Worksheets.Item("Alcohol Consumption")
$ws6.Range("A12").Formula = "= 7 +'Tobbaco USe'!A17"

# 7. Drugs and Marijuana Use: existing row 8, A8 = 3 + 'Alcohol Consumption'!A12
$ws7 = $wb.Worksheets.Item("Drugs and Marijuana Use")
$ws7.Range("A8").Formula = "= 3 +'Alcohol Consumption'!A12"

# 8. Immunization: new row 9, A9 = 4 + 'Drugs and Marijuana Use'!A8
$ws8 = $wb.Worksheets.Item("Immunization")
$ws8.Range("A9").Formula = "= 4 +'Drugs and Marijuana Use'!A8"

# 9. Urban Rural: new row 7, A7 = 2 + Immunization!A9
$ws9 = $wb.Worksheets.Item("Urban Rural")
$ws9.Range("A7").Formula = "= 2 +Immunization!A9"

# Update selections to match the post-edit cursor positions
[void]$ws1.Activate()
[void]$ws1.Range("A9").Select()

[void]$ws2.Activate()
[void]$ws2.Range("A13").Select()

[void]$ws3.Activate()
[void]$ws3.Range("A7").Select()

[void]$ws4.Activate()
[void]$ws4.Range("A11").Select()

[void]$ws5.Activate()
[void]$ws5.Range("A17").Select()

[void]$ws6.Activate()
[void]$ws6.Range("A12").Select()

[void]$ws7.Activate()
[void]$ws7.Range("A8").Select()

[void]$ws8.Activate()
[void]$ws8.Range("A9").Select()

[void]$ws9.Activate()
[void]$ws9.Range("A8").Select()

# The final active sheet is "Health Status" (tab index 3)
[void]$ws4.Activate()
